$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.002.75'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '2.305.62'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = "'303.56"
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = "'97.99"
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('E7').Value = '  -1.51%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = "'0.508"
$ws.Range('E9').Value = '  -1.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = "'36.00"
$ws.Range('E10').Value = '  +1.11%  '
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = "'18.16"
$ws.Range('E12').Value = '  +1.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = "'6.81"
$ws.Range('E14').Value = '  -0.95%  '
$ws.Range('D15').Value = '2.663.24'
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').Value = '2.299.60'
$ws.Range('E16').Value = '  -0.73%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = "'0.784"
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').Value = '42.937.10'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = "'12.64"
$ws.Range('E19').Value = '  -5.70%  '
$ws.Range('D20').Value = '0.0₃0905'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = "'6.04"
$ws.Range('E21').Value = '  -1.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = "'68.04"
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = "'236.78"
$ws.Range('E23').Value = '  -0.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = "'2.15"
$ws.Range('E24').Value = '  -1.87%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = "'2.47"
$ws.Range('E25').Value = '  +1.85%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = "'25.49"
$ws.Range('E28').Value = '  +3.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = "'165.40"
$ws.Range('E29').Value = '  -1.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = "'2.06"
$ws.Range('E30').Value = '  +0.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = "'9.08"
$ws.Range('E31').Value = '  -0.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = "'33.47"
$ws.Range('E32').Value = '  +1.27%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('E34').Value = '  +0.81%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = "'5.03"
$ws.Range('E35').Value = '  -2.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = "'17.04"
$ws.Range('E36').Value = '  -6.12%  '
$ws.Range('E37').Value = '  -1.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = "'0.0691"
$ws.Range('E38').Value = '  +0.35%  '
$ws.Range('E39').Value = '  -1.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = "'1.77"
$ws.Range('E40').Value = '  -1.04%  '
$ws.Range('E41').Value = '  -1.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = "'2.74"
$ws.Range('E42').Value = '  -0.54%  '
$ws.Range('D43').Value = '2.008.92'
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = "'0.0283"
$ws.Range('E44').Value = '  -1.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = "'10.08"
$ws.Range('E45').Value = '  -0.76%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = "'17.86"
$ws.Range('E46').Value = '  +2.07%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = "'2.09"
$ws.Range('E47').Value = '  -2.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = "'2.81"
$ws.Range('E48').Value = '  -0.56%  '
$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = "'2.88"
$ws.Range('E49').Value = '  +3.46%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = "'53.83"
$ws.Range('E50').Value = '  -1.06%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.532.34'
$ws.Range('E51').Value = '  +0.13%  '
